$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("I2").Value = 2.5
$ws.Range("Q2").Value = 3
$ws.Range("Y2").Value = 1.5
$ws.Range("AC2").Value = 3.5
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 0.142
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 0.142
$ws.Range("AE3").Value = 0.71
$ws.Range("E4").Value = 3.5
$ws.Range("I4").Value = 6.5
$ws.Range("O4").Value = 2
$ws.Range("P4").Value = 100
$ws.Range("Q4").Value = 3.5
$ws.Range("R4").Value = 0.142
$ws.Range("AC4").Value = 2.5
$ws.Range("AE4").Value = 0.781
$ws.Range("U5").Value = 1
$ws.Range("Y5").Value = 6
$ws.Range("E6").Value = 5.29
$ws.Range("I6").Value = 8.289999999999999
$ws.Range("M6").Value = 1.43
$ws.Range("Q6").Value = 10.14
$ws.Range("U6").Value = 1.57
$ws.Range("Y6").Value = 4.29
$ws.Range("AC6").Value = 7.14
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("Q7").Value = 7
$ws.Range("U7").Value = 1
$ws.Range("AC7").Value = 2
$ws.Range("AE7").Value = 0.71
$ws.Range("E9").Value = 4.33
$ws.Range("Q9").Value = 3.33
$ws.Range("W9").Value = 2
$ws.Range("X9").Value = 66.67
$ws.Range("Y9").Value = 1
$ws.Range("Z9").Value = 0.095
$ws.Range("AE9").Value = 0.71
$ws.Range("I10").Value = 3
